$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is forced to text format ("@") before assignment so that
# numeric-looking strings (e.g. "9.12", "0.999") are stored as literal
# text rather than being coerced to a floating point number, matching
# the inlineStr cell type used throughout this sheet. ClearFormats()
# afterwards removes the temporary number-format style so the cell
# keeps its original (unstyled) appearance.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '76.472.08'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.052.18'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.42%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '202.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.72'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.60%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E7').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.210'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +6.25%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.051.24'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.45%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.438'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('E11').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.15'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.12%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.615.73'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.60%  '
$ws.Range('E14').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.85%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.405.97'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000195'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.043.24'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.85%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.58'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.19%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.12'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.21%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.14'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('E21').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.37'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.60'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.207.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.38%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.44'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.85%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.97'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.60%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000112'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.40%  '
$ws.Range('E29').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.35'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.62%  '
$ws.Range('E31').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '507.78'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.96'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.03%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '20.86'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.04%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '163.13'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('E37').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.05%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '191.34'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.68%  '
$ws.Range('E40').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('E41').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.802'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +22.14%  '
$ws.Range('E44').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.13%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '42.31'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.79%  '
$ws.Range('E47').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.48'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.70%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.612'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.89%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.92'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.39%  '
$ws.Range('E51').ClearFormats()
